# Regenerate the localization-status report for a new handoff run.
# Old run id:  978f8624-989f-4489-a515-978a1512e981
# New run id:  3cbda7d4-8939-4a3b-9ac8-dbec18455361
# Old xliff content hash: 807481d03a69b79da67edfeece7f95a34130eb8c
# New xliff content hash: d46e06e4d1598a6f2baa80f96cce7e2e418f4fcb

$wb = $excel.ActiveWorkbook

$newId = "3cbda7d4-8939-4a3b-9ac8-dbec18455361"
$newHash = "d46e06e4d1598a6f2baa80f96cce7e2e418f4fcb"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("B2").Value = "e2e\$newId.md"
foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = "e2e\$newId.md"
}
$wsOverview.Range("G2").Value = "2016-08-21 21:07:50"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newId.md"
foreach ($hl in $wsZh.Hyperlinks) {
    $hl.TextToDisplay = "$newId.md"
}
$wsZh.Range("G2").Value = "$newId.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-21 21:07:46"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newId.md"
foreach ($hl in $wsDe.Hyperlinks) {
    $hl.TextToDisplay = "$newId.md"
}
$wsDe.Range("G2").Value = "$newId.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-21 21:07:50"
